$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 (the last existing header cell) onto the new H1 header
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New "Save" header column
$ws.Range("H1").Value = "Save"

# New Save value for the data row
$ws.Range("H2").Value = 1
